$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "四段文？第一段，弥勒请文殊回答大众疑惑",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "四段文？【答】第一段，弥勒请文殊回答大众疑惑",
    2
)
